$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Swap the F:V (match details) content between three row pairs.
#    (A and E columns -- index & match date -- stay where they are.)
# ---------------------------------------------------------------
function Swap-MatchRows($rowA, $rowB) {
    $va = $ws.Range("F$rowA`:V$rowA").Value()
    $vb = $ws.Range("F$rowB`:V$rowB").Value()
    $ws.Range("F$rowA`:V$rowA").Value = $vb
    $ws.Range("F$rowB`:V$rowB").Value = $va
}

Swap-MatchRows 284 285
Swap-MatchRows 297 298
Swap-MatchRows 346 347

# ---------------------------------------------------------------
# 2) Append the new match as row 349.
# ---------------------------------------------------------------

# Columns that are identical to the row above (348) can simply be
# copied (value + formatting) so their text type/format is preserved
# exactly (e.g. "2023" must stay text, not become numeric).
$ws.Cells.Item(348,2).Copy()
$ws.Cells.Item(349,2).PasteSpecial(-4104)   # B: pais

$ws.Cells.Item(348,3).Copy()
$ws.Cells.Item(349,3).PasteSpecial(-4104)   # C: torneio

$ws.Cells.Item(348,4).Copy()
$ws.Cells.Item(349,4).PasteSpecial(-4104)   # D: temporada

# A: Indice (new numeric value, but keep the bold/border style of A348)
$ws.Cells.Item(349,1).Value = 348
$ws.Cells.Item(348,1).Copy()
$ws.Cells.Item(349,1).PasteSpecial(-4122)   # xlPasteFormats

# E: data_partida (new date serial, keep the date-time style of E348)
$ws.Cells.Item(349,5).Value = 45235.91666666666
$ws.Cells.Item(348,5).Copy()
$ws.Cells.Item(349,5).PasteSpecial(-4122)   # xlPasteFormats

$ws.Cells.Item(349,6).Value = "Vitoria"
$ws.Cells.Item(349,7).Value = 1
$ws.Cells.Item(349,8).Value = "Vila Nova FC"
$ws.Cells.Item(349,9).Value = 1
$ws.Cells.Item(349,10).Value = 2.1
$ws.Cells.Item(349,11).Value = "30/10/2023 02:42"
$ws.Cells.Item(349,12).Value = 2.08
$ws.Cells.Item(349,13).Value = "05/11/2023 21:51"
$ws.Cells.Item(349,14).Value = 3.07
$ws.Cells.Item(349,15).Value = "30/10/2023 02:42"
$ws.Cells.Item(349,16).Value = 3.02
$ws.Cells.Item(349,17).Value = "05/11/2023 21:51"
$ws.Cells.Item(349,18).Value = 3.86
$ws.Cells.Item(349,19).Value = "30/10/2023 02:42"
$ws.Cells.Item(349,20).Value = 4.48
$ws.Cells.Item(349,21).Value = "05/11/2023 21:51"
$ws.Cells.Item(349,22).Value = "https://www.betexplorer.com/football/brazil/serie-b/vitoria-vila-nova-fc/dt90iwcs/"
